# Update cryptos list — price (D) and volume-change (E) columns, plus a
# ranking swap between Binance-Peg BSC-USD (was row 33) and EthereumClassic
# (was row 34), refreshed by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($Row, $D) {
    # Many price strings look like plain numbers ("32.50", "679.59", ...).
    # A bare .Value assignment lets Excel auto-coerce them to a Double and
    # silently drop the trailing zero / formatting, so force text storage,
    # write the value, then drop back to the default (unstyled) look so no
    # stray numeric format / quote-prefix style sticks to the cell.
    $cell = $ws.Cells.Item($Row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $D
    $cell.Style = "Normal"
}

function Set-Pct($Row, $E) {
    $ws.Cells.Item($Row, 5).Value = $E
}

Set-Price 2  "69.470.37"
Set-Pct   2  "  +0.07%  "

Set-Price 3  "3.690.45"
Set-Pct   3  "  -0.08%  "

Set-Pct   4  "  +0.08%  "

Set-Price 5  "679.59"
Set-Pct   5  "  -1.01%  "

Set-Price 6  "161.28"
Set-Pct   6  "  +0.07%  "

Set-Price 7  "0.999"
Set-Pct   7  "  +0.02%  "

Set-Pct   8  "  -0.04%  "

Set-Pct   9  "  -0.05%  "

Set-Price 10 "7.17"
Set-Pct   10 "  -0.93%  "

Set-Price 11 "0.440"
Set-Pct   11 "  +0.28%  "

Set-Pct   12 "  -0.41%  "

Set-Price 13 "4.312.31"
Set-Pct   13 "  -0.12%  "

Set-Price 14 "32.50"
Set-Pct   14 "  -0.53%  "

Set-Price 15 "3.687.81"
Set-Pct   15 "  +0.10%  "

Set-Price 16 "69.420.45"
Set-Pct   16 "  -0.08%  "

Set-Pct   17 "  +2.64%  "

Set-Price 18 "16.06"
Set-Pct   18 "  +0.69%  "

Set-Price 19 "6.49"
Set-Pct   19 "  +0.31%  "

Set-Price 20 "471.62"
Set-Pct   20 "  -0.72%  "

Set-Pct   21 "  -1.24%  "

Set-Pct   22 "  +0.33%  "

Set-Price 23 "80.41"
Set-Pct   23 "  +0.92%  "

Set-Price 24 "3.836.30"
Set-Pct   24 "  -0.01%  "

Set-Pct   25 "  -0.16%  "

Set-Pct   26 "  -0.04%  "

Set-Price 27 "10.90"
Set-Pct   27 "  -1.47%  "

Set-Pct   28 "  -1.12%  "

Set-Pct   29 "  -0.53%  "

Set-Pct   30 "  -1.46%  "

Set-Price 31 "2.01"
Set-Pct   31 "  -0.85%  "

Set-Price 32 "6.59"
Set-Pct   32 "  -1.70%  "

# Rows 33 and 34 swapped ranking order: EthereumClassic moved up to 33,
# Binance-Peg BSC-USD moved down to 34.
$ws.Cells.Item(33, 2).Value = "EthereumClassic"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-Price 33 "27.00"
Set-Pct   33 "  +0.74%  "

$ws.Cells.Item(34, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-Price 34 "0.998"
Set-Pct   34 "  -0.30%  "

Set-Price 35 "3.680.06"
Set-Pct   35 "  +0.35%  "

Set-Pct   36 "  +1.54%  "

Set-Price 37 "8.45"
Set-Pct   37 "  +2.54%  "

Set-Price 38 "6.21"
Set-Pct   38 "  +1.41%  "

Set-Pct   39 "  -0.01%  "

Set-Pct   40 "  -1.22%  "

Set-Pct   41 "  -0.03%  "

Set-Pct   42 "  -0.85%  "

Set-Price 43 "168.55"
Set-Pct   43 "  +1.31%  "

Set-Pct   44 "  -0.22%  "

Set-Price 45 "46.67"
Set-Pct   45 "  -2.64%  "

Set-Price 46 "2.74"
Set-Pct   46 "  -0.36%  "

Set-Price 47 "0.000281"
Set-Pct   47 "  +1.43%  "

Set-Price 48 "28.12"
Set-Pct   48 "  -1.53%  "

Set-Price 49 "1.29"
Set-Pct   49 "  -2.41%  "

Set-Pct   50 "  -3.16%  "

Set-Pct   51 "  +0.38%  "
